# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits previously-merged "<word><space>" text runs back into two
# separate runs: one holding the bare word and one holding a single
# space, mirroring the pre-consolidation OOXML (each `<a:t>` chunk ends
# up as its own `<a:r>` with an empty `<a:rPr/>`).
#
# `Characters(start, length).Text = sameText` on a PowerPoint TextRange
# forces that sub-span to materialize as its own run without touching
# the rest of the paragraph's text, so re-assigning every token back to
# itself (in place, left to right) reproduces the desired run split.

function Split-TextRuns {
    param(
        $TextRange,
        [string[]]$Tokens
    )

    $pos = 1
    foreach ($tok in $Tokens) {
        $len = $tok.Length
        if ($len -gt 0) {
            $TextRange.Characters($pos, $len).Text = $tok
        }
        $pos += $len
    }
}

$p = $ppt.ActivePresentation

# --- slide 1 : "Slide 1 (Content)" ---------------------------------
Split-TextRuns $p.Slides.Item(1).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "1", " ", "(Content)")

# --- slide 2 : "Slide 2 (Content)" ---------------------------------
Split-TextRuns $p.Slides.Item(2).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "2", " ", "(Content)")

# --- slide 3 : "Slide 3 (Content)" ---------------------------------
Split-TextRuns $p.Slides.Item(3).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "3", " ", "(Content)")

# --- slide 4 : "Slide 4 (Content)" ---------------------------------
Split-TextRuns $p.Slides.Item(4).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "4", " ", "(Content)")

# --- slide 5 : "Slide 5 (Two Content)" ------------------------------
Split-TextRuns $p.Slides.Item(5).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "5", " ", "(Two", " ", "Content)")

# --- slide 6 : "Slide 6 (Two Content Right)" + "an image" ----------
Split-TextRuns $p.Slides.Item(6).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "6", " ", "(Two", " ", "Content", " ", "Right)")
Split-TextRuns $p.Slides.Item(6).Shapes.Item("TextBox 3").TextFrame.TextRange `
    @("an", " ", "image")

# --- slide 7 : "Slide 7 (Content with Caption)" + "An image" -------
Split-TextRuns $p.Slides.Item(7).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "7", " ", "(Content", " ", "with", " ", "Caption)")
Split-TextRuns $p.Slides.Item(7).Shapes.Item("TextBox 3").TextFrame.TextRange `
    @("An", " ", "image")

# --- slide 8 : "Slide 8 (Comparison)" + "An image" ------------------
Split-TextRuns $p.Slides.Item(8).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "8", " ", "(Comparison)")
Split-TextRuns $p.Slides.Item(8).Shapes.Item("TextBox 3").TextFrame.TextRange `
    @("An", " ", "image")

# --- slide 9 : "Slide 10 (Content)" ---------------------------------
Split-TextRuns $p.Slides.Item(9).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "10", " ", "(Content)")

# --- slide 10 : "Slide 11 (Content)" --------------------------------
Split-TextRuns $p.Slides.Item(10).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "11", " ", "(Content)")

# --- slide 11 : "Slide 12 (Content)" --------------------------------
Split-TextRuns $p.Slides.Item(11).Shapes.Title.TextFrame.TextRange `
    @("Slide", " ", "12", " ", "(Content)")
